# Class2_WhatIsDataViz.pptx - "Added homework to the slide deck."
#
# 1. Slide 10 ("In class exercise"): drop the two trailing paragraphs that
#    told students to read the Quealy post / note the process (that content
#    moves to the new Homework slide), and let the placeholder's autofit
#    recompute now that there is less text.
# 2. A new slide 11 ("Homework") is appended at the end of the deck with the
#    reading / repo-update assignment.

$LCURLY = [char]0x201C
$RCURLY = [char]0x201D
$RSQUO  = [char]0x2019

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 10: trim the homework reminder off the in-class-exercise slide.
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$body10 = $s10.Shapes.Item(2)

$body10.TextFrame.TextRange.Text = (
    "Google " + $LCURLY + "Money on the Bench" + $RCURLY + "`r" +
    "What" + $RSQUO + "s the lead?`r" +
    "Nut graph?`r" +
    "Body?`r" +
    "Background and details?`r" +
    "What did you take away from this? "
)
# Fewer/shorter lines now fit without shrinking - drop the stale
# fontScale/lnSpcReduction normAutofit hints by letting autosize recompute.
$body10.TextFrame.AutoSize = 2

# ---------------------------------------------------------------------
# New slide 11: "Homework", same "Title and Content" layout as slide 10.
# ---------------------------------------------------------------------
$s11 = $p.Slides.Add($p.Slides.Count + 1, 2)

$title11 = $s11.Shapes.Item(1)
$title11.TextFrame.TextRange.Text = "Homework"

$body11 = $s11.Shapes.Item(2)
$tr = $body11.TextFrame.TextRange
$tr.Text = "Read "
$null = $tr.InsertAfter("Tufte")
$null = $tr.InsertAfter(" Chapter 1 Part 1 on Graphical Excellence (pg. 13-51).")

$null = $tr.InsertAfter("`rRead Kevin ")
$null = $tr.InsertAfter("Quealy" + $RSQUO + "s")
$null = $tr.InsertAfter(" post on ")
$null = $tr.InsertAfter("chartsnthings")
$null = $tr.InsertAfter(" about Money on the Bench (third link in your Google search for money on the bench). ")
$null = $tr.InsertAfter("Note the ")
$null = $tr.InsertAfter("process.")

$null = $tr.InsertAfter("`rUpdate your forked course repository. Hint: Google " + $LCURLY + "update forked repository" + $RCURLY + " ")

$null = $tr.InsertAfter("`r")

$body11.TextFrame.AutoSize = 2
